$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [double]"0.9247055945318959"
$ws.Range("E2").Value = [double]"0.9247055945318959"

$ws.Range("D3").Value = [double]"0.9999650201456074"
$ws.Range("E3").Value = [double]"0.9999650201456074"

$ws.Range("D4").Value = [double]"0.004155324731252437"
$ws.Range("E4").Value = [double]"0.004155324731252437"

$ws.Range("D5").Value = [double]"1.071617584048023E-32"
$ws.Range("E5").Value = [double]"1.071617584048023E-32"

$ws.Range("D6").Value = [double]"0.9067412681750768"
$ws.Range("E6").Value = [double]"0.9067412681750768"

$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = [double]"2.53727126121521"
